$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New consolidated values: each row combines the token name with a
# Python tuple/list-style repr string of its card text lines.
$values = @(
    '(''Beast'', [''Token Creature — Beast'', ''3/3''])',
    '(''Cat'', [''Token Creature — Cat'', ''2/2''])',
    '(''Drake'', [''Token Creature — Drake'', ''Flying'', ''2/2''])',
    '(''Goat'', [''Token Creature — Goat'', ''0/1''])',
    '(''Goblin'', [''Token Creature — Goblin'', ''1/1''])',
    '(''Hellion'', [''Token Creature — Hellion'', ''4/4''])',
    '(''Liliana of the Dark Realms Emblem'', [''Emblem — Liliana'', "Swamps you control have ‘{T}: Add {B}{B}{B}{B}.''"])',
    '(''Saproling'', [''Token Creature — Saproling'', ''1/1''])',
    '(''Soldier'', [''Token Creature — Soldier'', ''1/1''])',
    '(''Wurm'', [''Token Creature — Wurm'', ''6/6''])',
    '(''Zombie'', [''Token Creature — Zombie'', ''2/2''])'
)

# Clear out the old data rows (2 through 35) before writing the new,
# shorter list of rows (2 through 12).
$ws.Range("A2:A35").ClearContents()

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}
